# Reproduce the commit: "Made sheet2 the same as reader/sheet2 ..."
# The substantive edit is on "Sheet2 - Numbers": a new column AA (col 27) is
# populated with the values 100..129 for rows 1..30, the sheet becomes the
# active sheet/tab, and the selection moves to AA1:AA30 (with AA1 as the
# active cell). Sheet4's page setup paper size also changes from "Any"(0)
# to A4 (9).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")

# Fill the new column AA (column 27) with 100, 101, ..., 129 on rows 1-30.
for ($row = 1; $row -le 30; $row++) {
    $ws2.Cells.Item($row, 27).Value = 99 + $row
}

# Make Sheet2 the active sheet/tab and select the newly written range,
# matching the new <selection activeCell="AA1" sqref="AA1:AA30"/>.
$ws2.Activate()
$ws2.Range("AA1:AA30").Select()

# Sheet4 - Dates: pageSetup paperSize changes from 0 to 9 (A4).
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
